$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.86%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.82%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.864"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.26%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06421"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.58%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'1.31%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.190"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-7.54%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.97%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.55%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.06%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07499"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.36%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.64%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08973"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.94%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.13%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006360"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.15%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006109"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.25%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.473"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.304"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.53%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E21").Value = "'1.74%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.904"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.36%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1518"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'9.99%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04404"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001176"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.54%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.72%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'-1.74%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'1.63%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.04119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.68%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006821"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.98%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.29%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001890"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-18.25%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01166"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.69%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.24%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'9.30%"
$ws.Range("E46").Style = "Normal"
